# Verify Get Single User.xlsx — "add New TC,TS, and latest reports"
#
# Net effect of the change (5-row data table in Sheet1, rows 2-7):
#   * D2: "Missing API key"  -> "Missing API key." (trailing period added)
#   * A new test-case is inserted logically at row 3
#       A3="aaaa", B3=1, C3=401, D3="Invalid API key.", E3="Check authorization"
#     which pushes the former row-3 record ("reqres-free-v1", -, 200, -,
#     "get all user response") down into row 4.
#   * Rows 5-7 are untouched.
#
# Implemented as direct cell writes (values + number/text formatting) so the
# final grid matches exactly, rather than a literal Rows.Insert (which would
# grow the sheet by one row — not what the target file shows: row count
# stays at 7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteAll = -4104
$xlPasteFormats = -4122

function Copy-CellAllAndFormat($srcAddr, $dstAddr) {
    # Copies both the value/formula AND the cell style/number-format from
    # $srcAddr to $dstAddr. A single PasteSpecial(xlPasteAll) in this host
    # carries the value but not the style, so we do a second
    # PasteSpecial(xlPasteFormats) pass to bring the formatting along too.
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial($xlPasteAll)
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = $false
}

# --- Phase 1: pull out everything that still needs to read an *original*
#     (pre-edit) value/style from a cell before that cell gets overwritten.

# Old row-3 record slides down into row 4.
Copy-CellAllAndFormat "A3" "A4"
$ws.Range("C4").Value = 200

# E3/E4 swap places (new row3 needs old E4's text+style, new row4 needs old
# E3's text+style) - route the old E3 content through a scratch cell.
Copy-CellAllAndFormat "E3" "Z1"
Copy-CellAllAndFormat "E4" "E3"
Copy-CellAllAndFormat "Z1" "E4"
$ws.Range("Z1").Clear()

# New row3's D cell reuses old row4's style (s=1); text is replaced after.
Copy-CellAllAndFormat "D4" "D3"
$ws.Range("D3").Value = "Invalid API key."

# New row3's B cell is an exact copy of old row4's B cell (1, s=2).
Copy-CellAllAndFormat "B4" "B3"

# --- Phase 2: the donor cells in row 4 are now fully absorbed elsewhere;
#     clear them completely (value + style) so they disappear like the
#     target file shows (no <c> element at all for B4/D4).
$ws.Range("B4").Clear()
$ws.Range("D4").Clear()

# --- Phase 3: remaining straightforward value edits (existing style is
#     already correct on these cells, so a plain value assignment suffices).
$ws.Range("A3").Value = "aaaa"
$ws.Range("C3").Value = 401
$ws.Range("D2").Value = "Missing API key."
